# Slotted_Screws.xlsx: rework the per-sheet header block.
#
# Table_1: row 1 (previously the text header) becomes a numeric index row
# (0..11), row 2 (previously the lone "18-8 Stainless Steel" marker cell)
# becomes the real text header, and every data row's Material/Surface
# column (L) is stamped with "18-8 Stainless Steel".
#
# Table_2: a new row is inserted above row 1 so every data row shifts down
# by one. The new row 1 becomes the numeric index row (0..11, keeping the
# bold/bordered header style), and the row below it (the old text header)
# is rewritten to the same text-header template used on Table_1 (its
# thread_size/material_surface columns are dropped).

$wb = $excel.ActiveWorkbook

$headerLabels = @("Lg.", "Threading", "HeadDia.", "Head Ht.", "DriveSize", "TensileStrength, psi", "Specifications Met", "Pkg.Qty.", "", "Pkg.", "", "")

# ---------------------------------------------------------------------
# Table_1
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table_1")

# Row 2 gets the text header labels (was just "18-8 Stainless Steel" in A2).
for ($i = 0; $i -lt 12; $i++) {
    $ws1.Cells.Item(2, $i + 1).Value2 = $headerLabels[$i]
}

# Row 1 becomes the numeric index row 0..11 (keeps its existing bold style).
for ($i = 0; $i -lt 12; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value2 = $i
}

# Every data row (3..75) gets "18-8 Stainless Steel" in column L.
for ($r = 3; $r -le 75; $r++) {
    $ws1.Cells.Item($r, 12).Value2 = "18-8 Stainless Steel"
}

# ---------------------------------------------------------------------
# Table_2
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Table_2")

# Insert a blank row at the top; every existing row (header + data) shifts
# down by one, so the old header (row 1) becomes row 2 and the old data
# rows 2..28 become rows 3..29.
$ws2.Rows.Item(1).Insert()

# Move the bold/bordered header style from row 2 (the pushed-down old
# header) up onto the new row 1, then strip that formatting from row 2
# so it goes back to plain/default like every other data row.
$ws2.Range("A2:L2").Copy()
$ws2.Range("A1:L1").PasteSpecial(-4122)
$ws2.Range("A2:L2").ClearFormats()

# Row 1 becomes the numeric index row 0..11.
for ($i = 0; $i -lt 12; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value2 = $i
}

# Row 2 becomes the text header labels (thread_size/material_surface
# columns intentionally left blank, matching Table_1's header row).
for ($i = 0; $i -lt 12; $i++) {
    $ws2.Cells.Item(2, $i + 1).Value2 = $headerLabels[$i]
}
